$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44497
$ws.Range("J2").Value = 150
$ws.Range("K2").Value = 6000
$ws.Range("L2").Value = 6500
$ws.Range("M2").Value = 6333
$ws.Range("P2").Value = 253

# Row 3
$ws.Range("D3").Value = 44503
$ws.Range("J3").Value = 250
$ws.Range("K3").Value = 9000
$ws.Range("L3").Value = 10000
$ws.Range("M3").Value = 9400
$ws.Range("O3").Value = 'Provincia de Melipilla'
$ws.Range("P3").Value = 376

# Row 4
$ws.Range("D4").Value = 44517
$ws.Range("J4").Value = 130
$ws.Range("K4").Value = 6000
$ws.Range("L4").Value = 6500
$ws.Range("M4").Value = 6269
$ws.Range("O4").Value = 'Región Metropolitana'
$ws.Range("P4").Value = 251

# Row 5
$ws.Range("D5").Value = 44523
$ws.Range("J5").Value = 100
$ws.Range("M5").Value = 9500
$ws.Range("O5").Value = 'Región Metropolitana'
$ws.Range("P5").Value = 380

# Row 6
$ws.Range("D6").Value = 44476
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 7000
$ws.Range("L6").Value = 7500
$ws.Range("M6").Value = 7250
$ws.Range("P6").Value = 290

# Row 7
$ws.Range("D7").Value = 44545
$ws.Range("J7").Value = 140
$ws.Range("K7").Value = 14000
$ws.Range("L7").Value = 15000
$ws.Range("M7").Value = 14429
$ws.Range("O7").Value = 'Provincia de Chacabuco'
$ws.Range("P7").Value = 577

# Row 8
$ws.Range("D8").Value = 44505
$ws.Range("J8").Value = 180
$ws.Range("K8").Value = 6000
$ws.Range("L8").Value = 6500
$ws.Range("M8").Value = 6222
$ws.Range("O8").Value = 'Región del Maule'
$ws.Range("P8").Value = 249

# Row 9
$ws.Range("D9").Value = 44384
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 12000
$ws.Range("L9").Value = 13000
$ws.Range("M9").Value = 12500
$ws.Range("O9").Value = 'Región de Coquimbo'
$ws.Range("P9").Value = 500

# Row 10
$ws.Range("D10").Value = 44533
$ws.Range("J10").Value = 180
$ws.Range("K10").Value = 8000
$ws.Range("L10").Value = 8500
$ws.Range("M10").Value = 8222
$ws.Range("O10").Value = 'Región del Maule'
$ws.Range("P10").Value = 329

# Row 11
$ws.Range("D11").Value = 44526
$ws.Range("J11").Value = 100
$ws.Range("K11").Value = 7500
$ws.Range("L11").Value = 8000
$ws.Range("M11").Value = 7750
$ws.Range("P11").Value = 310

# Row 12
$ws.Range("D12").Value = 44467
$ws.Range("K12").Value = 8000
$ws.Range("L12").Value = 9000
$ws.Range("M12").Value = 8500
$ws.Range("P12").Value = 340

# Row 13
$ws.Range("D13").Value = 44540
$ws.Range("J13").Value = 140
$ws.Range("K13").Value = 11000
$ws.Range("L13").Value = 12000
$ws.Range("M13").Value = 11429
$ws.Range("O13").Value = 'Región del Maule'
$ws.Range("P13").Value = 457

# Row 14
$ws.Range("D14").Value = 44509
$ws.Range("J14").Value = 100
$ws.Range("K14").Value = 6500
$ws.Range("L14").Value = 7000
$ws.Range("M14").Value = 6750
$ws.Range("O14").Value = 'Región Metropolitana'
$ws.Range("P14").Value = 270

# Row 15
$ws.Range("D15").Value = 44316
$ws.Range("K15").Value = 16000
$ws.Range("L15").Value = 18000
$ws.Range("M15").Value = 17000
$ws.Range("P15").Value = 680

# Row 17
$ws.Range("D17").Value = 44482
$ws.Range("J17").Value = 430
$ws.Range("K17").Value = 8000
$ws.Range("L17").Value = 8500
$ws.Range("M17").Value = 8267
$ws.Range("O17").Value = 'Región de O''Higgins'
$ws.Range("P17").Value = 331

# Row 18
$ws.Range("D18").Value = 44188
$ws.Range("K18").Value = 18000
$ws.Range("L18").Value = 20000
$ws.Range("M18").Value = 19000
$ws.Range("P18").Value = 760

# Row 19
$ws.Range("D19").Value = 44351
$ws.Range("J19").Value = 100
$ws.Range("K19").Value = 15000
$ws.Range("L19").Value = 16000
$ws.Range("M19").Value = 15500
$ws.Range("O19").Value = 'Región Metropolitana'
$ws.Range("P19").Value = 620

# Row 21
$ws.Range("D21").Value = 44162

# Row 22
$ws.Range("D22").Value = 44537
$ws.Range("J22").Value = 160
$ws.Range("K22").Value = 8500
$ws.Range("L22").Value = 9000
$ws.Range("M22").Value = 8719
$ws.Range("P22").Value = 349

# Row 23
$ws.Range("D23").Value = 44483
$ws.Range("J23").Value = 350
$ws.Range("K23").Value = 5500
$ws.Range("L23").Value = 6000
$ws.Range("M23").Value = 5714
$ws.Range("P23").Value = 229

# Row 24
$ws.Range("D24").Value = 44498
$ws.Range("J24").Value = 220
$ws.Range("K24").Value = 7000
$ws.Range("L24").Value = 7500
$ws.Range("M24").Value = 7273
$ws.Range("O24").Value = 'Región Metropolitana'
$ws.Range("P24").Value = 291

# Row 25
$ws.Range("D25").Value = 44335
$ws.Range("K25").Value = 18000
$ws.Range("L25").Value = 20000
$ws.Range("M25").Value = 19000
$ws.Range("O25").Value = 'Provincia de Limarí'
$ws.Range("P25").Value = 760

# Row 26
$ws.Range("D26").Value = 44160
$ws.Range("K26").Value = 9000
$ws.Range("L26").Value = 10000
$ws.Range("M26").Value = 9500
$ws.Range("P26").Value = 380
